$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D32").Value = "Spark DataFrame vs Pandas DataFrame"
$ws.Range("E32").Value = "https://dodonam.tistory.com/434"

$ws.Range("D36").Value = "Noisy Label Learning"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/409"

$ws.Range("D46").Value = "[특허청] 2023년 6월, 생물정보학(Bioinformatics 채용), 바이오 인포매틱스 분야 특허*실용신안 심사"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/543"

$ws.Range("D51").Value = "NAS (network attached storage)는 네트워크에 연결된 저장소"
$ws.Range("E51").Value = "https://bskyvision.com/entry/NAS-network-attached-storage%EB%8A%94-%EB%84%A4%ED%8A%B8%EC%9B%8C%ED%81%AC%EC%97%90-%EC%97%B0%EA%B2%B0%EB%90%9C-%EC%A0%80%EC%9E%A5%EC%86%8C"
